# Update "想去人数" (interest count) values in the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 359
$ws1.Range("F3").Value = 759
$ws1.Range("F4").Value = 261
$ws1.Range("F5").Value = 762
$ws1.Range("F6").Value = 1866

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 359
$ws4.Range("F3").Value = 759
$ws4.Range("F4").Value = 261
$ws4.Range("F7").Value = 762
$ws4.Range("F8").Value = 1866
